$d = $word.ActiveDocument

# 1. Update the title/heading paragraph text.
$d.Content.Find.Execute(
    "Implante de Marcapasso Fisiológico com Monitoramento Remoto (Azure™)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Marcapasso Fisiológico (Azure His/CSP – monitoramento remoto)", 2)

# 2. Remove the subtitle paragraph entirely
#    ("Estimulação hisiana/septal com telemonitorização.")
$d.Paragraphs(2).Range.Delete()

# 3. Update each material list item: add a "• " bullet prefix and
#    simplify/rewrite the wording.
$d.Content.Find.Execute(
    "Gerador – Azure™",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Gerador Azure – marcapasso dupla câmara com monitoramento remoto", 2)

$d.Content.Find.Execute(
    "Bainha His – C315™",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Bainha C315 para His/septo", 2)

$d.Content.Find.Execute(
    "Eletrodo His / Septal – 3830",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Eletrodo 3830", 2)

$d.Content.Find.Execute(
    "Eletrodo Atrial – 5076-52",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Eletrodo 5076-52", 2)

$d.Content.Find.Execute(
    "Ferramenta de Corte",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Ferramenta de corte", 2)

$d.Content.Find.Execute(
    "Fio Guia",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Fio guia", 2)

$d.Content.Find.Execute(
    "Introdutor – 2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Introdutor – 2", 2)
